# "Changed title to 12 pt font"
#
# The built-in "Title" paragraph style currently carries an explicit
# 18pt (w:sz / w:szCs = 36 half-points) size override in its run
# properties. Update it to 12pt (24 half-points), for both the
# regular and complex-script/bidi sizes, so the Title style matches
# the document's base 12pt body text.

$d = $word.ActiveDocument

$titleStyle = $d.Styles("Title")
$titleStyle.Font.Size = 12
$titleStyle.Font.SizeBi = 12

# Relocate the "_GoBack" bookmark (Word's marker for "the place I was
# last editing") from the reference_docx code sample -- where it sat
# next to the word "template" -- onto the date line ("April 8, 2016"),
# which is the text actually touched by this edit.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $dateParagraph = $p
    }
}

if ($dateParagraph -ne $null) {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()

    $dateRange = $d.Range($dateParagraph.Range.Start, $dateParagraph.Range.End - 1)
    $d.Bookmarks.Add("_GoBack", $dateRange)
}
